# Separate compiling data code from analysis:
# - "Raw" sheet: set F15:F33 all to 0.5, and C27:C33 from 150 to 125
# - "Data" sheet: set E16:E34 all to 0.5, and B28:B34 from 150 to 125
# Also update the active selections on those two sheets to match the
# authored workbook (cosmetic, but reflects where the author was last
# working while separating the sheets).

$wb = $excel.ActiveWorkbook

# --- "Raw" worksheet ---
$wsRaw = $wb.Worksheets.Item("Raw")
$wsRaw.Range("F15:F33").Value = 0.5
$wsRaw.Range("C27:C33").Value = 125
$wsRaw.Activate()
$wsRaw.Range("F6:F33").Select()

# --- "Data" worksheet ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("E16:E34").Value = 0.5
$wsData.Range("B28:B34").Value = 125
$wsData.Activate()
$wsData.Range("D7").Select()
